# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F4").Value  = 128
$wsExhibit.Range("F12").Value = 588
$wsExhibit.Range("F20").Value = 254
$wsExhibit.Range("F24").Value = 6442
$wsExhibit.Range("F29").Value = 120
$wsExhibit.Range("F32").Value = 1250
$wsExhibit.Range("F38").Value = 226

# 全部类型 sheet updates (same events, different row positions)
$wsAll.Range("F17").Value = 588
$wsAll.Range("F25").Value = 254
$wsAll.Range("F31").Value = 6442
$wsAll.Range("F36").Value = 1250
$wsAll.Range("F45").Value = 226
